# Update the GLM model-summary table.
# The three per-group models (PE and PVC / PS and Tile / All Data) were
# re-fit: the "no3" term used in model 1 (PE and PVC) and model 3
# (All Data) was swapped out for "time", and the refreshed
# estimate / std.error / statistic / p.value numbers for the intercept and
# that term are written in for each of those two models.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Model 1: "PE and PVC" (rows 2-4: intercept, time, nh4) ---------------
$ws.Range("B3").Value = "time"
$ws.Range("B4").Value = "nh4"

$ws.Range("C2").Value = 19.65
$ws.Range("D2").Value = 0.738
$ws.Range("E2").Value = 26.63
$ws.Range("F2").Value = [double]"3.29E-41"

$ws.Range("C3").Value = -0.1957
$ws.Range("D3").Value = 0.007594
$ws.Range("E3").Value = -25.76
$ws.Range("F3").Value = [double]"3.401E-40"

$ws.Range("C4").Value = -0.1438
$ws.Range("D4").Value = 0.005487
$ws.Range("E4").Value = -26.21
$ws.Range("F4").Value = [double]"1.012E-40"

# --- Model 3: "All Data" (rows 10-12: intercept, time, nh4) ---------------
$ws.Range("B11").Value = "time"
$ws.Range("B12").Value = "nh4"

$ws.Range("C10").Value = 13.07
$ws.Range("D10").Value = 1.925
$ws.Range("E10").Value = 6.789
$ws.Range("F10").Value = [double]"9.129E-11"

$ws.Range("C11").Value = -0.1333
$ws.Range("D11").Value = 0.01973
$ws.Range("E11").Value = -6.758
$ws.Range("F11").Value = [double]"1.091E-10"

$ws.Range("C12").Value = -0.09854
$ws.Range("D12").Value = 0.01431
$ws.Range("E12").Value = -6.887
$ws.Range("F12").Value = [double]"5.172E-11"
